# Applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.108.22'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.560.34'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.24%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  -3.13%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.10'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0590'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0864'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.781.25'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.549.45'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.517'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.120.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '214.80'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.22%  '
$ws.Range('E19').Value = '  -1.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.74%  '
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('E26').Value = '  -2.58%  '
$ws.Range('E27').Value = '  -1.85%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  -1.66%  '
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('E32').Value = '  -2.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.383.45'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('E34').Value = '  -1.37%  '
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('E36').Value = '  -1.70%  '
$ws.Range('E37').Value = '  -3.11%  '
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.810'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.514'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.52%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.987'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E45').Value = '  -0.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.26'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.694.72'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.37'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.94%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₇0982'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.27%  '
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('E51').Value = '  +0.10%  '
